$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark from the "Sprawdzić parametry..." paragraph.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# 2. Swap the text of the two list items:
#    "Sprawdzić działanie zaznaczania ostatniego elementu , czy przycisk się pojawia"
#    <->
#    "Sprawdzić zależności miedzy checkboxami"
$oldA = "Sprawdzić działanie zaznaczania ostatniego elementu , czy przycisk się pojawia"
$oldB = "Sprawdzić zależności miedzy checkboxami"
$placeholder = "##SWAP_PLACEHOLDER##"

$d.Content.Find.Execute($oldA, $true, $false, $false, $false, $false, $true, 1, $false, $placeholder, 2)
$d.Content.Find.Execute($oldB, $true, $false, $false, $false, $false, $true, 1, $false, $oldA, 2)
$d.Content.Find.Execute($placeholder, $true, $false, $false, $false, $false, $true, 1, $false, $oldB, 2)

# 3. Add a new "_GoBack" bookmark spanning the paragraph that now reads
#    "Sprawdzić zależności miedzy checkboxami" (text only, not the mark).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith($oldB)) {
        $r = $d.Range($p.Range.Start, $p.Range.End - 1)
        $d.Bookmarks.Add("_GoBack", $r)
        break
    }
}
